$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New Fitness (column C) values, keyed by Generation (column B) range.
# Generation 0..33  -> 7734
# Generation 34..45 -> 7343
# Generation 46..250 -> 7293
for ($gen = 0; $gen -le 250; $gen++) {
    $row = $gen + 2
    if ($gen -le 33) {
        $value = 7734
    } elseif ($gen -le 45) {
        $value = 7343
    } else {
        $value = 7293
    }
    $ws.Cells.Item($row, 3).Value = $value
}
